$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("D2").Value = '92.513.61'
$ws.Range("E2").Value = '  -1.90%  '
$ws.Range("B3").Value = 'Ethereum'
$ws.Range("D3").Value = '3.332.95'
$ws.Range("E3").Value = '  -3.36%  '
$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("B5").Value = 'Solana'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '230.85'
$ws.Range("E5").Value = '  -2.71%  '
$ws.Range("B6").Value = 'BNB'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '613.90'
$ws.Range("E6").Value = '  -4.60%  '
$ws.Range("B7").Value = 'XRP'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.39'
$ws.Range("E7").Value = '  -3.29%  '
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.388'
$ws.Range("E8").Value = '  -3.95%  '
$ws.Range("B9").Value = 'USDC'
$ws.Range("E9").Value = '  -0.02%  '
$ws.Range("B10").Value = 'Cardano'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.955'
$ws.Range("E10").Value = '  -1.76%  '
$ws.Range("B11").Value = 'LidoStakedEther'
$ws.Range("D11").Value = '3.333.30'
$ws.Range("E11").Value = '  -3.29%  '
$ws.Range("B12").Value = 'Avalanche'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.76'
$ws.Range("E12").Value = '  +1.44%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.195'
$ws.Range("E13").Value = '  -1.34%  '
$ws.Range("B14").Value = 'Toncoin'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.22'
$ws.Range("E14").Value = '  +1.21%  '
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("D15").Value = '92.286.86'
$ws.Range("E15").Value = '  -1.89%  '
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("D16").Value = '3.947.60'
$ws.Range("E16").Value = '  -3.52%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("E17").Value = '  -3.06%  '
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.01'
$ws.Range("E18").Value = '  -4.60%  '
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("D19").Value = '3.330.82'
$ws.Range("E19").Value = '  -3.46%  '
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.25'
$ws.Range("E20").Value = '  -2.10%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.96'
$ws.Range("E21").Value = '  -4.37%  '
$ws.Range("B22").Value = 'SuiNetwork'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.40'
$ws.Range("E22").Value = '  +4.78%  '
$ws.Range("B23").Value = 'BitcoinCash'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '491.99'
$ws.Range("E23").Value = '  -1.51%  '
$ws.Range("B24").Value = 'Stellar'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.433'
$ws.Range("E24").Value = '  -13.11%  '
$ws.Range("B25").Value = 'NEARProtocol'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.51'
$ws.Range("E25").Value = '  -1.50%  '
$ws.Range("B26").Value = 'PEPE'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000182'
$ws.Range("E26").Value = '  -5.51%  '
$ws.Range("B27").Value = 'Litecoin'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '90.25'
$ws.Range("E27").Value = '  -3.93%  '
$ws.Range("B28").Value = 'Aptos'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.95'
$ws.Range("E28").Value = '  -1.06%  '
$ws.Range("B29").Value = 'WrappedeETH'
$ws.Range("D29").Value = '3.499.65'
$ws.Range("E29").Value = '  -3.58%  '
$ws.Range("B30").Value = 'Dai'
$ws.Range("E30").Value = '  -0.02%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.09'
$ws.Range("E31").Value = '  -6.39%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.137'
$ws.Range("E32").Value = '  -1.12%  '
$ws.Range("B33").Value = 'PancakeSwap'
$ws.Range("E33").Value = '  -5.26%  '
$ws.Range("B34").Value = 'Binance-PegBSC-USD'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  +0.20%  '
$ws.Range("B35").Value = 'Cronos'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.172'
$ws.Range("E35").Value = '  -2.82%  '
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '28.21'
$ws.Range("E36").Value = '  -5.77%  '
$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '564.56'
$ws.Range("E37").Value = '  +1.95%  '
$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.524'
$ws.Range("E38").Value = '  -5.58%  '
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("E39").Value = '  -3.02%  '
$ws.Range("B40").Value = 'USDe'
$ws.Range("E40").Value = '  -0.05%  '
$ws.Range("B41").Value = 'Kaspa'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.148'
$ws.Range("E41").Value = '  -2.34%  '
$ws.Range("B42").Value = 'Fetch.AI'
$ws.Range("E42").Value = '  -4.86%  '
$ws.Range("B43").Value = 'ARBITRUM'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.861'
$ws.Range("E43").Value = '  -5.39%  '
$ws.Range("B44").Value = 'WhiteBITCoin'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '23.73'
$ws.Range("E44").Value = '  -1.29%  '
$ws.Range("B45").Value = 'MantraDAO'
$ws.Range("C45").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.65'
$ws.Range("E45").Value = '  +1.07%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0410'
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("B47").Value = 'ImmutableX'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.66'
$ws.Range("E47").Value = '  -3.40%  '
$ws.Range("B48").Value = 'Filecoin'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '5.40'
$ws.Range("E48").Value = '  -3.00%  '
$ws.Range("B49").Value = 'Stacks'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.10'
$ws.Range("E49").Value = '  -3.92%  '
$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '52.23'
$ws.Range("E50").Value = '  -1.91%  '
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.92'
$ws.Range("E51").Value = '  -1.60%  '
